$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 25, pushing the "Vignette" (was row25) and
#     "Testing of sumby" (was row26) rows down to 26 and 27 respectively ---
$ws.Rows(25).Insert()

# --- Row 24 (the gap left behind by the insert) becomes a new bug entry ---
$ws.Range("A24").Value = "SB"
$ws.Range("B24").Value = 43411
$ws.Range("C24").Value = "Not clear how to change the root calling programme name from the default of ""main.R"", when needed."

# give B24 the same date number-format used throughout column B/D (style carried from B23)
$ws.Range("B23").Copy()
$ws.Range("B24").PasteSpecial(-4122)

# --- Row 13: note a new "can't recreate" comment in column D ---
$ws.Range("D13").Value = "14/11/18 can't recreate this bug"

# --- Row 25 is now blank content-wise, but keeps the date style in column B ---
$ws.Range("B23").Copy()
$ws.Range("B25").PasteSpecial(-4122)

# --- Row 27 ("Testing of sumby") gains a resolved-date in column D ---
$ws.Range("D27").Value = 43418
$ws.Range("B23").Copy()
$ws.Range("D27").PasteSpecial(-4122)

# --- restore the active selection shown in the workbook ---
$ws.Range("D28").Select()
